# Corrected Thingiverse and GitHub links
#
# The "Vendor URL" column (K) for the hardware rows (15-18) pointed at the
# wrong Thingiverse/GitHub links (the STDHT_MGN9_NOVA variant's links,
# copy-pasted from a sibling BOM workbook). Point them at the correct
# thing/repo path for this part (CARRIAGE/STDHT_MGN12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLink = "https://www.thingiverse.com/thing:4810748`nalternate: https://github.com/MirageC79/HextrudORT/tree/main/files/CARRIAGE/STDHT_MGN12/STL"

$ws.Range("K15").Value = $newLink
$ws.Range("K16").Value = $newLink
$ws.Range("K17").Value = $newLink
$ws.Range("K18").Value = $newLink

# Restore the view/selection state as it was left after the edit: scrolled
# down so D15 is the top-left visible cell, with E17 as the active cell.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("E17").Select()
